$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New trade row (row 31, trade #30) appended below the existing data.
    $ws.Cells.Item(31, 1).Value = 30
    # Leading "'" forces text so the date-looking string isn't auto-converted
    # to a date serial number (matches how Excel stores a typed date string
    # verbatim when the column has no date formatting). Resetting the style
    # back to "Normal" afterwards drops the quote-prefix formatting flag that
    # the apostrophe trick leaves behind, so the cell ends up as plain text
    # with the default (unstyled) format, same as the rest of the sheet.
    $ws.Cells.Item(31, 2).Value = "'2026-02-17"
    $ws.Cells.Item(31, 2).Style = "Normal"
    $ws.Cells.Item(31, 3).Value = "08:23:09"
    $ws.Cells.Item(31, 4).Value = "MarketMaking"
    $ws.Cells.Item(31, 5).Value = "UP"
    $ws.Cells.Item(31, 6).Value = 0.55
    # Exit Price is blank for this still-open trade.
    $ws.Cells.Item(31, 7).Value = "'"
    $ws.Cells.Item(31, 7).Style = "Normal"
    $ws.Cells.Item(31, 8).Value = "OPEN"
    $ws.Cells.Item(31, 9).Value = 0
    $ws.Cells.Item(31, 10).Value = 0
    $ws.Cells.Item(31, 11).Value = 98.9686731447009
    $ws.Cells.Item(31, 12).Value = 0
    $ws.Cells.Item(31, 13).Value = 0
    $ws.Cells.Item(31, 14).Value = 0.6
    $ws.Cells.Item(31, 15).Value = "Normal spread capture: 19600 bps"
    # Exit Reason is blank for this still-open trade.
    $ws.Cells.Item(31, 16).Value = "'"
    $ws.Cells.Item(31, 16).Style = "Normal"
    $ws.Cells.Item(31, 17).Value = 0
}
